$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section header row 69: "Final Calculator " (bold section header, like row 65) ---
$ws.Range("B65").Copy()
$ws.Range("B69").PasteSpecial(-4122)
$ws.Range("B69").Value = "Final Calculator "

# --- Data rows 70-76 ---
# Copy the formatting template from an existing similar data row (66) so the
# new rows inherit the same per-column styles (A: date format, B/C/E: default).
# Copy column-by-column (rather than A:E as one block) so empty D cells in
# the template don't get materialised in the new rows.
$ws.Range("A66").Copy()
$ws.Range("A70:A76").PasteSpecial(-4122)
$ws.Range("B66").Copy()
$ws.Range("B70:B76").PasteSpecial(-4122)
$ws.Range("C66").Copy()
$ws.Range("C70:C76").PasteSpecial(-4122)
$ws.Range("E66").Copy()
$ws.Range("E70:E76").PasteSpecial(-4122)

# Dates for all new rows
$ws.Range("A70").Value = "11/16/2020"
$ws.Range("A71").Value = "11/16/2020"
$ws.Range("A72").Value = "11/16/2020"
$ws.Range("A73").Value = "11/16/2020"
$ws.Range("A74").Value = "11/16/2020"
$ws.Range("A75").Value = "11/16/2020"
$ws.Range("A76").Value = "11/16/2020"

# "Passed" column (all "Y") for new rows
$ws.Range("C70").Value = "Y"
$ws.Range("C71").Value = "Y"
$ws.Range("C72").Value = "Y"
$ws.Range("C73").Value = "Y"
$ws.Range("C74").Value = "Y"
$ws.Range("C75").Value = "Y"
$ws.Range("C76").Value = "Y"

# Description / Notes text -- written in the same order the shared-string
# table grew in the source file: B71, E71, B70, E70, E72, B73, B72, E73,
# B74, E74, B75, E75, B76, E76.
$ws.Range("B71").Value = "Verifying variables resulting as integers"
$ws.Range("E71").Value = "Used the var_dump tool to display the array."

$ws.Range("B70").Value = "Testing cloud function on server"
$ws.Range("E70").Value = "Used the function-add function to display test results"

$ws.Range("E72").Value = "Positive numbers calculate as expected"

$ws.Range("B73").Value = "Testing negative numbers for add tests"
$ws.Range("B72").Value = "Testing positive numbers  for add tests"
$ws.Range("E73").Value = "Results and negative numbers return as expected"

$ws.Range("B74").Value = "Two zero's test"
$ws.Range("E74").Value = "Returns back a zero"

$ws.Range("B75").Value = "Negative zero test"
$ws.Range("E75").Value = "Does not return a -0 as an output."

$ws.Range("B76").Value = "Numbers only test"
$ws.Range("E76").Value = "Cannot type letters or non-negative/positive symbols, with the exception of  e (euler) (?)"

# --- View state: selection follows the new bottom of data ---
$ws.Range("E76").Select()
